$d = $word.ActiveDocument

# 1. Remove the standalone "Meta description: ..." paragraph that sits
#    right after the title (Heading1) paragraph near the top of the doc.
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2. At the end of the document, the paragraph that used to hold the
#    "Please create an image..." image-prompt text gets replaced by two
#    paragraphs: a new bold title line, followed by the meta-description
#    text (now italicised) that used to live near the top of the doc.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

# Target only the paragraph's content, not its trailing paragraph mark,
# so InsertXML's "replace the range" semantics swap the text/runs in
# place instead of leaving a stray empty paragraph behind.
$target = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)

$replacementXml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r/>
            <w:r>
              <w:rPr><w:b/></w:rPr>
              <w:t>Play Crystal Quest: Deep Jungle for Free | Slot Review</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r/>
            <w:r>
              <w:rPr><w:i/></w:rPr>
              <w:t>Explore the rainforest and enjoy the unlimited multiplier and bonus round of free spins in Crystal Quest: Deep Jungle. Play for free and read our review.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$target.InsertXML($replacementXml)
